$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.057.36"
$ws.Range("E2").Value = "  -3.09%  "

$ws.Range("D3").Value = "3.002.43"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'533.22"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("D6").Value = "'133.83"
$ws.Range("E6").Value = "  +0.66%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "2.997.20"
$ws.Range("E8").Value = "  -2.02%  "

$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  -3.41%  "

$ws.Range("D11").Value = "'6.09"
$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("D14").Value = "'34.20"
$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "3.489.05"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "61.043.90"
$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").Value = "3.003.61"
$ws.Range("E18").Value = "  -2.17%  "

$ws.Range("D19").Value = "'6.60"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'463.93"
$ws.Range("E20").Value = "  -3.57%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "'0.674"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("D23").Value = "'6.94"
$ws.Range("E23").Value = "  -1.94%  "

$ws.Range("D24").Value = "'79.24"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").Value = "'12.02"
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").Value = "'7.85"
$ws.Range("E28").Value = "  -2.66%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").Value = "'25.50"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  +3.51%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'5.46"
$ws.Range("E33").Value = "  +2.51%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'55.35"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("E35").Value = "  -3.45%  "

$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  -1.70%  "

$ws.Range("D37").Value = "'455.30"
$ws.Range("E37").Value = "  -5.01%  "

$ws.Range("D38").Value = "3.216.96"
$ws.Range("E38").Value = "  +3.92%  "

$ws.Range("D39").Value = "'0.0785"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "'0.0383"
$ws.Range("E40").Value = "  -2.55%  "

$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").Value = "'8.15"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").Value = "'27.71"
$ws.Range("E43").Value = "  +14.02%  "

$ws.Range("E44").Value = "  -5.82%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "'0.246"
$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").Value = "'119.22"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("D50").Value = "0.0₃0494"
$ws.Range("E50").Value = "  -7.99%  "

$ws.Range("E51").Value = "  +7.75%  "
